$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) to the right of the existing "sum" column (G).
# Copy G1's formatting (bold font, borders, centered alignment header style)
# onto H1 so it reuses the same cell style as the other header cells.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"

# Fill the new column's data rows with 0 (unstyled, like the other numeric columns)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
